# 0.0.14 - Add Nodal and Member Concentrated Load
#
# The "type" column (C) for the Node Loads (rows 3-5) changes from the
# generic "L" marker to "C" (Concentrated), and the now-unused trailing
# blank row (row 9) is removed. Selection moves to C10 to reflect the new
# bottom of the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "C"
$ws.Range("C4").Value = "C"
$ws.Range("C5").Value = "C"

$ws.Rows.Item(9).Delete()

$ws.Range("C10").Select()
